$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 114: new data row for 2025-10-27, station "四方坪站充电量(kw)"
$ws.Range("A114").Value = 45957
$ws.Range("B114").Value = "四方坪站充电量(kw)"
$ws.Range("C114").Value = 529.37099999999987
$ws.Range("D114").Value = 1180.3180000000002
$ws.Range("E114").Value = 370.83
$ws.Range("F114").Value = 510.38800000000003
$ws.Range("G114").Value = 376.40499999999997
$ws.Range("H114").Value = 555.9380000000001
$ws.Range("I114").Value = 395.12200000000001
$ws.Range("J114").Value = 273.95600000000002
$ws.Range("K114").Value = 70.944000000000003
$ws.Range("L114").Value = 118.11599999999999
$ws.Range("M114").Value = 209.08799999999999
$ws.Range("N114").Value = 201.95999999999998
$ws.Range("O114").Value = 800.4259999999997
$ws.Range("P114").Value = 954.38800000000015
$ws.Range("Q114").Value = 584.83400000000017
$ws.Range("R114").Value = 307.79000000000002
$ws.Range("S114").Value = 187.97
$ws.Range("T114").Value = 200.46
$ws.Range("U114").Value = 81.64
$ws.Range("V114").Value = 187.08
$ws.Range("W114").Value = 18.549999999999997
$ws.Range("X114").Value = 13.6
$ws.Range("Y114").Value = 29.93
$ws.Range("Z114").Value = 93.722999999999999

# Row 115: new data row for 2025-10-27, station "高岭站充电量(kw)"
$ws.Range("A115").Value = 45957
$ws.Range("B115").Value = "高岭站充电量(kw)"
$ws.Range("C115").Value = 425.51800000000003
$ws.Range("D115").Value = 269.279
$ws.Range("E115").Value = 107.187
$ws.Range("F115").Value = 99.769000000000005
$ws.Range("G115").Value = 65.304000000000002
$ws.Range("H115").Value = 82.179000000000002
$ws.Range("I115").Value = 133.63200000000001
$ws.Range("J115").Value = 124.744
$ws.Range("K115").Value = 151.97499999999999
$ws.Range("L115").Value = 158.934
$ws.Range("M115").Value = 88.349000000000004
$ws.Range("N115").Value = 188.28500000000003
$ws.Range("O115").Value = 454.41200000000003
$ws.Range("P115").Value = 415.29500000000002
$ws.Range("Q115").Value = 271.20599999999996
$ws.Range("R115").Value = 428.75600000000003
$ws.Range("S115").Value = 369.78
$ws.Range("T115").Value = 99.925999999999988
$ws.Range("U115").Value = 101.342
$ws.Range("V115").Value = 20.172000000000001
$ws.Range("W115").Value = 52.210999999999999
$ws.Range("X115").Value = 40.272999999999996
$ws.Range("Y115").Value = 82.632999999999996
$ws.Range("Z115").Value = 7.641

# Update the active-cell selection to match the post-edit state
[void]$ws.Range("H121").Select()
